# Update cryptos list values (price and 1h volume change) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.910.26'
$ws.Range("E2").Value = '  -0.28%  '
$ws.Range("D3").Value = '1.551.69'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  -0.56%  '
$ws.Range("D5").Value = '206.59'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = '0.487'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").Value = '22.01'
$ws.Range("E8").Value = '  +1.79%  '
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").Value = '0.0595'
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").Value = '1.772.60'
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").Value = '1.543.30'
$ws.Range("E13").Value = '  -0.85%  '
$ws.Range("D14").Value = '3.75'
$ws.Range("E14").Value = '  +0.91%  '
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '26.899.00'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = '61.57'
$ws.Range("E17").Value = '  -0.57%  '
$ws.Range("E18").Value = '  +2.88%  '
$ws.Range("D19").Value = '217.14'
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("D20").Value = '7.29'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").Value = '153.43'
$ws.Range("E25").Value = '  +0.67%  '
$ws.Range("D26").Value = '6.64'
$ws.Range("E26").Value = '  -0.29%  '
$ws.Range("D27").Value = '14.97'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("E28").Value = '  +0.69%  '
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("E30").Value = '  +1.35%  '
$ws.Range("D31").Value = '1.08'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("E33").Value = '  +3.75%  '
$ws.Range("D34").Value = '1.411.47'
$ws.Range("E34").Value = '  +0.64%  '
$ws.Range("E35").Value = '  +2.30%  '
$ws.Range("D36").Value = '0.974'
$ws.Range("E36").Value = '  +2.23%  '
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("D38").Value = '0.0166'
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").Value = '0.526'
$ws.Range("E39").Value = '  +0.69%  '
$ws.Range("D40").Value = '0.806'
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("E41").Value = '  -0.57%  '
$ws.Range("E42").Value = '  +2.76%  '
$ws.Range("D43").Value = '2.30'
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").Value = '64.43'
$ws.Range("E46").Value = '  -0.96%  '
$ws.Range("D47").Value = '1.686.29'
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").Value = '87.19'
$ws.Range("E48").Value = '  +1.03%  '
$ws.Range("E49").Value = '  +1.60%  '
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +2.61%  '
$ws.Range("E51").Value = '  +0.26%  '
